$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Name:" paragraph and the blank paragraph that follows it
#    (paragraphs 2 and 3), so the title is immediately followed by the
#    "Ann has a 5% chance..." paragraph.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(2)
$pEnd   = $d.Paragraphs.Item(3)
$d.Range($pStart.Range.Start, $pEnd.Range.End).Delete()

# ---------------------------------------------------------------------------
# 2. Strip the trailing " (5 pts)" from each of the three numbered
#    questions (leaving a single trailing space) and mark each question's
#    paragraph with <w:contextualSpacing w:val="0"/>. InsertXML replaces the
#    whole paragraph (pPr + run) in one shot so both edits land atomically.
#    We locate each paragraph by scanning (not Find) because a Find hit's
#    Range covers only the matched text, not the trailing paragraph mark,
#    and InsertXML-ing a whole <w:p> onto that short range duplicates text.
# ---------------------------------------------------------------------------
function Get-ParagraphIndexByPrefix($doc, [string]$prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Set-QuestionParagraph($doc, [string]$prefix, [string]$newText) {
    $idx = Get-ParagraphIndexByPrefix $doc $prefix
    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
        '<w:contextualSpacing w:val="0"/></w:pPr><w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $doc.Paragraphs.Item($idx).Range.InsertXML($xml) | Out-Null
}

Set-QuestionParagraph $d "Calculate the actuarially fair annual premium" `
    "Calculate the actuarially fair annual premium for insuring Ann on her own. Do the same for Betty and Clara. "

Set-QuestionParagraph $d "If Ann, Betty, and Clara were all in the same insurance pool" `
    "If Ann, Betty, and Clara were all in the same insurance pool and the insurance company was forced to charge each of them the same premium, what premium would they have to charge to each to break even?  "

Set-QuestionParagraph $d "If the insurance company were to charge the premiums from problem 2" `
    "If the insurance company were to charge the premiums from problem 2, what would happen next?  "

# ---------------------------------------------------------------------------
# 3. Remove the blank "buffer" paragraphs that used to give students room to
#    answer under each question (bottom-up so earlier paragraph indices stay
#    valid as later ranges are deleted).
# ---------------------------------------------------------------------------
$d.Range($d.Paragraphs.Item(17).Range.Start, $d.Paragraphs.Item(20).Range.End).Delete()
$d.Range($d.Paragraphs.Item(11).Range.Start, $d.Paragraphs.Item(15).Range.End).Delete()
$d.Range($d.Paragraphs.Item(6).Range.Start, $d.Paragraphs.Item(9).Range.End).Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
